$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-5 (columns B..L), recalculated results
$data = @{
    2 = @(0.1666666666666667, 0.08058607483302643, 0.08608092516697356, 306.66746069682, 307.0230317067672, 307.9332303460307, 306.66746069682, 308.9183798995069, 309.8179771970586, 325.6796778876795, 0.7790629194950628)
    3 = @(0.2166666666666667, 0.1055058401668167, 0.1111611598331833, 310.4785776821093, 310.7534595537888, 311.4438362575175, 310.4785776821093, 312.2190145692873, 312.8986012346134, 325.0958879058764, 0.6527110314384649)
    4 = @(0.2666666666666667, 0.1304563450453853, 0.1362106549546147, 312.8574984796292, 313.0816395885642, 313.6375222630601, 312.8574984796292, 314.2767951394444, 314.8223685259341, 324.7310275843784, 0.4986528269812863)
    5 = @(0.3166666666666667, 0.1554241421537281, 0.1612428578462719, 314.4840917792948, 314.6733446756259, 315.1385109103733, 314.4840917792948, 315.6825332862946, 316.1380456860362, 324.4813890718348, 0.316888550462615)
}

$cols = @("B","C","D","E","F","G","H","I","J","K","L")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# New row 6
$ws.Range("A6").Value = 4
$row6 = @(0.3666666666666666, 0.1804026295365585, 0.1862643704634415, 315.6665988481321, 315.8303792725541, 316.2302458203413, 315.6665988481321, 316.7037701672701, 317.0946492976096, 324.299834870174, 0.1074183015888076)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])6").Value = $row6[$i]
}

# Apply same style as A2:A5 to A6 (header-like bold/centered/bordered style)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
